$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.952.75"
$ws.Range("E2").Value = "'  +0.34%  "
$ws.Range("D3").Value = "'1.817.65"
$ws.Range("E3").Value = "'  +0.42%  "
$ws.Range("E4").Value = "'  +0.16%  "
$ws.Range("D5").Value = "'309.86"
$ws.Range("E5").Value = "'  +0.04%  "
$ws.Range("E6").Value = "'  +0.15%  "
$ws.Range("D7").Value = "'0.4677"
$ws.Range("E7").Value = "'  +0.59%  "
$ws.Range("D8").Value = "'0.3665"
$ws.Range("E8").Value = "'  -0.82%  "
$ws.Range("D9").Value = "'0.07343"
$ws.Range("E9").Value = "'  -0.30%  "
$ws.Range("D10").Value = "'0.8722"
$ws.Range("E10").Value = "'  -0.12%  "
$ws.Range("D11").Value = "'20.24"
$ws.Range("E11").Value = "'  -1.10%  "
$ws.Range("D12").Value = "'1.819.82"
$ws.Range("E12").Value = "'  -1.76%  "
$ws.Range("D13").Value = "'5.405"
$ws.Range("E13").Value = "'  +0.84%  "
$ws.Range("E14").Value = "'  +0.79%  "
$ws.Range("D15").Value = "'6.507"
$ws.Range("E15").Value = "'  +0.01%  "
$ws.Range("D16").Value = "'91.38"
$ws.Range("E16").Value = "'  -0.09%  "
$ws.Range("E17").Value = "'  +0.20%  "
$ws.Range("D18").Value = "'0.000008709"
$ws.Range("E18").Value = "'  +0.11%  "
$ws.Range("E19").Value = "'  +0.07%  "
$ws.Range("D20").Value = "'14.65"
$ws.Range("E20").Value = "'  -0.48%  "
$ws.Range("D21").Value = "'26.968.93"
$ws.Range("E21").Value = "'  +0.29%  "
$ws.Range("D22").Value = "'5.292"
$ws.Range("E22").Value = "'  -0.43%  "
$ws.Range("D23").Value = "'10.59"
$ws.Range("E23").Value = "'  +0.69%  "
$ws.Range("D24").Value = "'2.043.60"
$ws.Range("E24").Value = "'  -0.43%  "
$ws.Range("E25").Value = "'  -0.57%  "
$ws.Range("D26").Value = "'150.87"
$ws.Range("D27").Value = "'18.36"
$ws.Range("E27").Value = "'  -0.30%  "
$ws.Range("D28").Value = "'2.151"
$ws.Range("E28").Value = "'  +0.20%  "
$ws.Range("D29").Value = "'5.256"
$ws.Range("E29").Value = "'  -0.87%  "
$ws.Range("D30").Value = "'117.03"
$ws.Range("E30").Value = "'  +1.10%  "
$ws.Range("D31").Value = "'0.08898"
$ws.Range("E31").Value = "'  +0.12%  "
$ws.Range("D32").Value = "'0.7584"
$ws.Range("E32").Value = "'  +0.72%  "
$ws.Range("D33").Value = "'1.161"
$ws.Range("E33").Value = "'  +0.78%  "
$ws.Range("D34").Value = "'4.496"
$ws.Range("E35").Value = "'  -0.14%  "
$ws.Range("E36").Value = "'  +0.14%  "
$ws.Range("D37").Value = "'1.090"
$ws.Range("E37").Value = "'  -0.86%  "
$ws.Range("D38").Value = "'0.05290"
$ws.Range("E38").Value = "'  +0.65%  "
$ws.Range("D39").Value = "'0.01945"
$ws.Range("E39").Value = "'  -0.93%  "
$ws.Range("D40").Value = "'2.969"
$ws.Range("E40").Value = "'  +1.62%  "
$ws.Range("D43").Value = "'0.5289"
$ws.Range("E43").Value = "'  -0.07%  "
$ws.Range("D44").Value = "'0.1652"
$ws.Range("E44").Value = "'  -0.60%  "
$ws.Range("D45").Value = "'8.428"
$ws.Range("E45").Value = "'  -0.18%  "
$ws.Range("D46").Value = "'0.4867"
$ws.Range("E46").Value = "'  -1.41%  "
$ws.Range("D47").Value = "'10.49"
$ws.Range("E47").Value = "'  +1.82%  "
$ws.Range("E48").Value = "'  +0.16%  "
$ws.Range("D49").Value = "'103.58"
$ws.Range("E49").Value = "'  +0.45%  "
$ws.Range("D50").Value = "'1.664"
$ws.Range("E50").Value = "'  -0.48%  "
$ws.Range("E51").Value = "'  +0.18%  "

# Row 41/42 swap: RenderToken moves to row 41, FraxShare moves to row 42
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "'2.373"
$ws.Range("E41").Value = "'  -2.26%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "'7.164"
$ws.Range("E42").Value = "'  -0.12%  "
